$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.458.62'
$ws.Range('E2').Value = '  +1.52%  '
$ws.Range('D3').Value = '1.676.58'
$ws.Range('E3').Value = '  +2.28%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''218.95'
$ws.Range('E5').Value = '  +2.33%  '
$ws.Range('D6').Value = '''0.5335'
$ws.Range('E6').Value = '  +1.79%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +4.23%  '
$ws.Range('D9').Value = '''0.06407'
$ws.Range('E9').Value = '  +1.75%  '
$ws.Range('D10').Value = '''21.85'
$ws.Range('E10').Value = '  +6.04%  '
$ws.Range('D11').Value = '''0.07793'
$ws.Range('E11').Value = '  +1.49%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.682.08'
$ws.Range('E12').Value = '  +2.48%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''4.514'
$ws.Range('E13').Value = '  +2.57%  '
$ws.Range('D14').Value = '''0.5587'
$ws.Range('E14').Value = '  +1.25%  '
$ws.Range('D15').Value = '0.0₅8329'
$ws.Range('E15').Value = '  +1.76%  '
$ws.Range('D16').Value = '''65.68'
$ws.Range('E16').Value = '  +1.17%  '
$ws.Range('D17').Value = '26.498.08'
$ws.Range('E17').Value = '  +1.72%  '
$ws.Range('D19').Value = '''4.797'
$ws.Range('E19').Value = '  +2.42%  '
$ws.Range('D20').Value = '''193.44'
$ws.Range('E20').Value = '  +2.65%  '
$ws.Range('E21').Value = '  +1.16%  '
$ws.Range('D22').Value = '''6.321'
$ws.Range('E22').Value = '  +2.63%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = '''0.1283'
$ws.Range('E24').Value = '  +6.32%  '
$ws.Range('D25').Value = '''140.51'
$ws.Range('E25').Value = '  -3.31%  '
$ws.Range('D26').Value = '''7.407'
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').Value = '''16.27'
$ws.Range('E27').Value = '  +3.11%  '
$ws.Range('D28').Value = '''1.443'
$ws.Range('E28').Value = '  +4.28%  '
$ws.Range('D29').Value = '''0.06287'
$ws.Range('E29').Value = '  +5.48%  '
$ws.Range('D30').Value = '''1.289'
$ws.Range('E30').Value = '  +2.77%  '
$ws.Range('D31').Value = '''3.608'
$ws.Range('E31').Value = '  +5.28%  '
$ws.Range('E32').Value = '  +1.87%  '
$ws.Range('D33').Value = '''1.697'
$ws.Range('E33').Value = '  +3.11%  '
$ws.Range('E34').Value = '  +3.06%  '
$ws.Range('D35').Value = '''0.6144'
$ws.Range('E35').Value = '  +9.28%  '
$ws.Range('E36').Value = '  +1.34%  '
$ws.Range('E37').Value = '  +0.92%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '''0.01631'
$ws.Range('E38').Value = '  +1.06%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '''6.138'
$ws.Range('E39').Value = '  +7.93%  '
$ws.Range('D40').Value = '1.095.82'
$ws.Range('E40').Value = '  +6.88%  '
$ws.Range('D41').Value = '''0.8652'
$ws.Range('E41').Value = '  +2.14%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = '''100.57'
$ws.Range('E43').Value = '  +0.47%  '
$ws.Range('D44').Value = '1.821.99'
$ws.Range('E44').Value = '  +1.97%  '
$ws.Range('E45').Value = '  +7.32%  '
$ws.Range('D46').Value = '''58.54'
$ws.Range('E46').Value = '  +5.08%  '
$ws.Range('D47').Value = '''8.165'
$ws.Range('E47').Value = '  +1.90%  '
$ws.Range('D48').Value = '''1.001'
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('D49').Value = '''0.05204'
$ws.Range('E49').Value = '  +1.26%  '
$ws.Range('E50').Value = '  +6.92%  '
$ws.Range('D51').Value = '''6.029'
$ws.Range('E51').Value = '  +2.20%  '
